# Applies the "Automatic update of files." edit:
#  - Column C (Förändrad) dates for rows 2-27 move from 2024-01-19 (45310)
#    to 2024-01-20 (45311).
#  - Row 27 gains the explicit default row height that all the other
#    existing data rows already carry.
#  - A brand-new record (row 28) is appended for case "A 2414-2024".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Bump the "Förändrad" date for the existing records (rows 2-27).
$ws.Range("C2:C27").Value = 45311

# 2) Row 27 now stores an explicit row height, like every other existing row.
$ws.Rows.Item(27).RowHeight = 15

# 3) Append the new row 28 record.
$ws.Range("A28").Value = "A 2414-2024"
$ws.Range("B28").Value = 45310
$ws.Range("C28").Value = 45311
$ws.Range("D28").Value = "OKÄNT"
$ws.Range("E28").Value = "OKÄNT"
$ws.Range("G28").Value = 7.2
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = 0
$ws.Range("N28").Value = 0
$ws.Range("O28").Value = 0
$ws.Range("P28").Value = 0
$ws.Range("Q28").Value = 0

# Match the formatting used by the rest of the table: dates in B/C are
# formatted as dates, and R has a wrap-text style (even though it is
# empty for this row, just like for the preceding rows).
$ws.Range("B28").NumberFormat = $ws.Range("B27").NumberFormat
$ws.Range("C28").NumberFormat = $ws.Range("C27").NumberFormat
$ws.Range("R28").WrapText = $true
